$d = $word.ActiveDocument

# 1. Update first paragraph text: "Coding in java" -> new intro line
$d.Paragraphs.Item(1).Range.Text = "Implement GUI of Take Quiz java files already coded."

# 2. Update second paragraph text: old "Take Quiz" note -> new note
$d.Paragraphs.Item(2).Range.Text = "Need to display current score, feedback, questions/answers"

# 3. Remove the old quiz-planning paragraphs ("- Make a quiz:" ... "Make that into GUI")
#    in their entirety, keeping only the trailing paragraph that holds the bookmark.
$pStart = $d.Paragraphs.Item(3)
$pEnd = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$midRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$midRange.Delete()

# 4. Clear the text of the final paragraph ("Keep on simple java programs ...")
#    but keep the now-empty paragraph mark so the _GoBack bookmark survives.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$textRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$textRange.Delete()
